# Update "想去人数" (F column) values on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

# Row -> new value mapping for column F (same updates apply to both sheets).
$updates = @{
    3  = 10879
    5  = 985
    6  = 199
    8  = 8340
    9  = 47
    10 = 471
    12 = 224
    13 = 139
    14 = 3333
    16 = 331
    17 = 40
    18 = 833
    19 = 135
    21 = 289
    22 = 132
    23 = 1855
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
